$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new cell values (homework grades) as described by the diff
$ws.Range("E4").Value = 5
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 5
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 5
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 5
$ws.Range("E11").Value = 5

# Update the active selection to E4 as shown in the diff
$ws.Range("E4").Select()
